$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new row of data (row 64) ---------------------------------
# A: Date, B: Time, C: Privacy, D: Post text, E: Reference (kept as text),
# F-M: reaction counts.
$ws.Range("A64").Value = 44161
$ws.Range("B64").Value = 0.9770833333333333
$ws.Range("C64").Value = "Friends"
$ws.Range("D64").Value = "DENSE FOG (Night Edition)"

# The reference value looks numeric but must stay text (like the other
# "Reference" column entries), so force the cell to text format first.
$ws.Range("E64").NumberFormat = "@"
$ws.Range("E64").Value = "10107853002117719"

$ws.Range("F64").Value = 3
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 1
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 1

# --- Refresh the worksheet's remembered sort state so it covers the new row
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A64"))
$sort.SetRange($ws.Range("A2:O64"))
$sort.Header = 0
$sort.Apply()

# --- Update the view: scroll position & active selection ------------------
$ws.Range("A37").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q67").Select()

# --- Nudge the workbook window position (best effort) ----------------------
try {
    $excel.ActiveWindow.Top = 37620
} catch {
}
